$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits on row 30 ---
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 5
$ws.Range("M30").Value = 5
$ws.Range("O30").Value = "ИЗМ"

# --- View changes: scroll frozen pane and move active selection ---
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("O31").Select()
